$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 0.999984930691919
$ws.Range("H2").Value = 0.777777777777778
$ws.Range("J2").Value = 3.625
$ws.Range("K2").Value = -1.13709905660377
$ws.Range("L2").Value = -1.57877571170818
$ws.Range("M2").Value = -0.762119266944709
$ws.Range("N2").Value = -31.3682498373455
$ws.Range("P2").Value = "Virtually certain improving"

# Row 3
$ws.Range("F3").Value = 0.306677177977373
$ws.Range("H3").Value = 0.884615384615385
$ws.Range("J3").Value = 11.13
$ws.Range("K3").Value = -0.0201334190003544
$ws.Range("L3").Value = -0.0600749620883886
$ws.Range("M3").Value = 0.0420333342810997
$ws.Range("N3").Value = -0.180893252473983
$ws.Range("P3").Value = "Unlikely increasing"

# Row 4
$ws.Range("F4").Value = 0.512609752963165
$ws.Range("H4").Value = 0.480769230769231
$ws.Range("L4").Value = -0.0013984096184896
$ws.Range("M4").Value = 0.0014237960936476

# Row 5
$ws.Range("H5").Value = 0.0204081632653061
$ws.Range("I5").Value = 1

# Row 6
$ws.Range("F6").Value = 0.999637525025523
$ws.Range("G6").Value = 0.846153846153846
$ws.Range("P6").Value = "Virtually certain improving"

# Row 7
$ws.Range("E7").Value = "ok"
$ws.Range("F7").Value = 0.0268465872455513
$ws.Range("G7").Value = 0.307692307692308
$ws.Range("H7").Value = 0.384615384615385
$ws.Range("J7").Value = 0.013
$ws.Range("K7").Value = 0.0010612712749192
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0.0026475603740742
$ws.Range("N7").Value = 8.16362519168642
$ws.Range("P7").Value = "Extremely unlikely improving"

# Row 8
$ws.Range("F8").Value = 0.025540585524228
$ws.Range("H8").Value = 0.711538461538462
$ws.Range("J8").Value = 7.685
$ws.Range("K8").Value = -0.0416505056890014
$ws.Range("L8").Value = -0.0728173008584146
$ws.Range("M8").Value = -0.0056405106263166
$ws.Range("N8").Value = -0.541971446831508
$ws.Range("P8").Value = "Extremely unlikely increasing"

# Row 9
$ws.Range("F9").Value = 0.0226481180001614
$ws.Range("H9").Value = 0.653846153846154
$ws.Range("J9").Value = 0.0165
$ws.Range("K9").Value = 0.000927811719668
$ws.Range("L9").Value = 0.0001495117257217
$ws.Range("M9").Value = 0.0016739661144839
$ws.Range("N9").Value = 5.62310133132144
$ws.Range("P9").Value = "Extremely unlikely improving"

# Row 10
$ws.Range("D10").Value = $false
$ws.Range("F10").Value = 0.805450278015326
$ws.Range("H10").Value = 0.646017699115044
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = -0.0720533586070691
$ws.Range("L10").Value = -0.235876896519715
$ws.Range("M10").Value = 0.0562613851946069
$ws.Range("N10").Value = -2.40177862023564
$ws.Range("P10").Value = "Likely improving"

# Row 11
$ws.Range("F11").Value = 0.747462760131149
$ws.Range("H11").Value = 0.776699029126214
$ws.Range("J11").Value = 11.08
$ws.Range("K11").Value = 0.0091528313569987
$ws.Range("L11").Value = -0.0164762198605148
$ws.Range("M11").Value = 0.0449850873621731
$ws.Range("N11").Value = 0.0826067812003498
$ws.Range("P11").Value = "Likely increasing"

# Row 12
$ws.Range("D12").Value = $false
$ws.Range("F12").Value = 0.996444972146462
$ws.Range("H12").Value = 0.285714285714286
$ws.Range("K12").Value = -0.0007249421898143
$ws.Range("L12").Value = -0.0011875806760913
$ws.Range("M12").Value = -0.0002722886734306
$ws.Range("N12").Value = -2.33852319294961
$ws.Range("P12").Value = "Virtually certain improving"

# Row 13
$ws.Range("G13").Value = 0.961165048543689
$ws.Range("H13").Value = 0.087378640776699

# Row 14
$ws.Range("F14").Value = 0.989159233238536
$ws.Range("G14").Value = 0.678571428571429
$ws.Range("H14").Value = 0.0714285714285714
$ws.Range("P14").Value = "Extremely likely improving"

# Row 15
$ws.Range("F15").Value = 0.0020193005112973
$ws.Range("G15").Value = 0.366071428571429
$ws.Range("H15").Value = 0.348214285714286
$ws.Range("J15").Value = 0.0097
$ws.Range("K15").Value = 0.0004155999976436
$ws.Range("M15").Value = 0.0010036347019392
$ws.Range("N15").Value = 4.28453605818184
$ws.Range("P15").Value = "Exceptionally unlikely improving"

# Row 16
$ws.Range("F16").Value = 0.110913500783029
$ws.Range("H16").Value = 0.555555555555556
$ws.Range("J16").Value = 7.71
$ws.Range("K16").Value = -0.0101663288547838
$ws.Range("L16").Value = -0.0220556696318059
$ws.Range("M16").Value = 0.0035038544255037
$ws.Range("N16").Value = -0.131858999413538

# Row 17
$ws.Range("F17").Value = 0.0006092593328019
$ws.Range("H17").Value = 0.526785714285714
$ws.Range("J17").Value = 0.01525
$ws.Range("K17").Value = 0.0008867171848115
$ws.Range("L17").Value = 0.000501717032967
$ws.Range("M17").Value = 0.0012490614936116
$ws.Range("N17").Value = 5.81453891679711
